$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "64.382.95"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.150.86"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "610.77"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "

# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "143.81"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -1.78%  "

# Row 7
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "3.150.28"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "

# Row 9
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.524"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +0.11%  "

# Row 10
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.150"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -0.54%  "

# Row 11
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "5.38"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -2.21%  "

# Row 12
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.470"
$r.Style = "Normal"
$ws.Range("E12").Value = "  -0.79%  "

# Row 13
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.0000255"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "

# Row 14
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "35.42"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -1.60%  "

# Row 15
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "3.668.94"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +0.45%  "

# Row 16
$ws.Range("E16").Value = "  +3.11%  "

# Row 17
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "64.396.46"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "3.152.79"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +0.46%  "

# Row 19
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "6.84"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -1.14%  "

# Row 20
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "476.21"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -0.65%  "

# Row 21
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "14.62"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +0.66%  "

# Row 22
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.725"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +3.22%  "

# Row 23
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "7.84"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +1.67%  "

# Row 24
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "13.69"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "

# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "84.34"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "

# Row 26
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "2.79"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -2.97%  "

# Row 28
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "8.51"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

# Row 29
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "7.32"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +6.99%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "0.118"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "2.10"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -4.64%  "

# Row 32
$ws.Range("E32").Value = "  +0.35%  "

# Row 33
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "26.42"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "

# Row 34
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "2.64"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -4.59%  "

# Row 35
$ws.Range("E35").Value = "  +1.52%  "

# Row 36
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "5.94"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "

# Row 37
$ws.Range("E37").Value = "  -2.64%  "

# Row 38
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.0₃0746"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +2.09%  "

# Row 39
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "3.08"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +3.87%  "

# Row 40
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "452.40"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "

# Row 41
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.0395"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -0.43%  "

# Row 42
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.118"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -1.46%  "

# Row 43
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "8.30"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "

# Row 44
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "2.841.53"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.265"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -1.30%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "2.26"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "

# Row 47
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "2.45"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +6.15%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "26.40"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "

# Row 49
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "

# Row 50
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.114"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "

# Row 51
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "34.24"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +3.53%  "
